# Auto-generated Excel COM-interop script to apply the Goblin Profits market-price refresh
# across all profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 8398.799999999999
$ws.Range("I13").Value = 7999.5
$ws.Range("J13").Value = 8665
$ws.Range("K13").Value = 7999.5
$ws.Range("L13").Value = 8665
$ws.Range("M13").Value = -7830.5
$ws.Range("N13").Value = -9003
$ws.Range("H19").Value = 1413.2727
$ws.Range("I19").Value = 889.4
$ws.Range("J19").Value = 1849.8334
$ws.Range("K19").Value = 889.4
$ws.Range("L19").Value = 1849.8334
$ws.Range("M19").Value = -714.4
$ws.Range("N19").Value = -2199.8334
$ws.Range("H62").Value = 62886.145
$ws.Range("I62").Value = 91156.336
$ws.Range("K62").Value = 91156.336
$ws.Range("M62").Value = -90532.336
$ws.Range("H65").Value = 62886.145
$ws.Range("I65").Value = 91156.336
$ws.Range("K65").Value = 455781.68
$ws.Range("M65").Value = -452661.68
$ws.Range("H80").Value = 1512.4
$ws.Range("J80").Value = 1815.25
$ws.Range("L80").Value = 5445.75
$ws.Range("N80").Value = -7441.75
$ws.Range("H83").Value = 1512.4
$ws.Range("J83").Value = 1815.25
$ws.Range("L83").Value = 16337.25
$ws.Range("N83").Value = -26321.25
$ws.Range("H127").Value = 2239.3635
$ws.Range("I127").Value = 1454.125
$ws.Range("J127").Value = 4333.3335
$ws.Range("K127").Value = 4362.375
$ws.Range("L127").Value = 13000.0005
$ws.Range("M127").Value = 597.625
$ws.Range("N127").Value = -22920.0005
$ws.Range("H135").Value = 858.45
$ws.Range("I135").Value = 858.45
$ws.Range("K135").Value = 7726.05
$ws.Range("M135").Value = -5191.05
$ws.Range("H137").Value = 1945.8948
$ws.Range("I137").Value = 2024
$ws.Range("K137").Value = 6072
$ws.Range("M137").Value = -3522
$ws.Range("H138").Value = 2055.845
$ws.Range("I138").Value = 837.4737
$ws.Range("J138").Value = 3458.818
$ws.Range("K138").Value = 2512.4211
$ws.Range("L138").Value = 10376.454
$ws.Range("M138").Value = 2627.5789
$ws.Range("N138").Value = -20656.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 333.33334
$ws.Range("I4").Value = 333.33334
$ws.Range("K4").Value = 333.33334
$ws.Range("M4").Value = -217.33334
$ws.Range("H45").Value = 80001470
$ws.Range("I45").Value = 50000776
$ws.Range("J45").Value = 90910820
$ws.Range("K45").Value = 50000776
$ws.Range("L45").Value = 90910820
$ws.Range("M45").Value = -50000399
$ws.Range("N45").Value = -90911574
$ws.Range("H61").Value = 5253
$ws.Range("I61").Value = 6595.091
$ws.Range("K61").Value = 6595.091
$ws.Range("M61").Value = -6383.091
$ws.Range("H63").Value = 3773
$ws.Range("J63").Value = 1994.5
$ws.Range("L63").Value = 1994.5
$ws.Range("N63").Value = -3366.5
$ws.Range("H66").Value = 3773
$ws.Range("J66").Value = 1994.5
$ws.Range("L66").Value = 9972.5
$ws.Range("N66").Value = -16836.5
$ws.Range("H74").Value = 2767.6667
$ws.Range("I74").Value = 3015.7144
$ws.Range("K74").Value = 3015.7144
$ws.Range("M74").Value = -2141.7144
$ws.Range("H77").Value = 2767.6667
$ws.Range("I77").Value = 3015.7144
$ws.Range("K77").Value = 15078.572
$ws.Range("M77").Value = -10710.572
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H132").Value = 1930.6086
$ws.Range("I132").Value = 1889.5
$ws.Range("K132").Value = 5668.5
$ws.Range("M132").Value = -3138.5
$ws.Range("H136").Value = 5253
$ws.Range("I136").Value = 6595.091
$ws.Range("K136").Value = 19785.273
$ws.Range("M136").Value = -17235.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3484.625
$ws.Range("I94").Value = 3163.842
$ws.Range("K94").Value = 3163.842
$ws.Range("M94").Value = -2712.842
$ws.Range("H99").Value = 2768.5625
$ws.Range("I99").Value = 1829.8
$ws.Range("J99").Value = 4333.1665
$ws.Range("K99").Value = 1829.8
$ws.Range("L99").Value = 4333.1665
$ws.Range("M99").Value = -331.8
$ws.Range("N99").Value = -7329.1665
$ws.Range("H134").Value = 4639.9
$ws.Range("I134").Value = 4914.143
$ws.Range("K134").Value = 14742.429
$ws.Range("M134").Value = -12207.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2036.6875
$ws.Range("I58").Value = 2072.4666
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 2072.4666
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -1869.4666
$ws.Range("N58").Value = -1906
$ws.Range("H62").Value = 12499.6875
$ws.Range("I62").Value = 9272.272000000001
$ws.Range("K62").Value = 9272.272000000001
$ws.Range("M62").Value = -8648.272000000001
$ws.Range("H65").Value = 12499.6875
$ws.Range("I65").Value = 9272.272000000001
$ws.Range("K65").Value = 46361.36
$ws.Range("M65").Value = -43241.36
$ws.Range("H132").Value = 3790.3333
$ws.Range("I132").Value = 3895.7058
$ws.Range("K132").Value = 11687.1174
$ws.Range("M132").Value = -9157.117400000001
$ws.Range("H134").Value = 3324.2273
$ws.Range("I134").Value = 3216.6
$ws.Range("J134").Value = 4400.5
$ws.Range("K134").Value = 9649.799999999999
$ws.Range("L134").Value = 13201.5
$ws.Range("M134").Value = -7114.799999999999
$ws.Range("N134").Value = -18271.5
$ws.Range("H136").Value = 2036.6875
$ws.Range("I136").Value = 2072.4666
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 6217.399800000001
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -3667.399800000001
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 7215298
$ws.Range("I17").Value = 20000294
$ws.Range("J17").Value = 112522.22
$ws.Range("K17").Value = 60000882
$ws.Range("L17").Value = 337566.66
$ws.Range("M17").Value = -60000713
$ws.Range("N17").Value = -337904.66
$ws.Range("H25").Value = 357.66666
$ws.Range("I25").Value = 311.5
$ws.Range("J25").Value = 450
$ws.Range("K25").Value = 934.5
$ws.Range("L25").Value = 1350
$ws.Range("M25").Value = -765.5
$ws.Range("N25").Value = -1688
$ws.Range("H30").Value = 357.66666
$ws.Range("I30").Value = 311.5
$ws.Range("J30").Value = 450
$ws.Range("K30").Value = 934.5
$ws.Range("L30").Value = 1350
$ws.Range("M30").Value = -832.5
$ws.Range("N30").Value = -1554
$ws.Range("H55").Value = 720586.1
$ws.Range("J55").Value = 7350.5
$ws.Range("L55").Value = 22051.5
$ws.Range("N55").Value = -22405.5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H94").Value = 8127.273
$ws.Range("I94").Value = 5700
$ws.Range("J94").Value = 8666.666999999999
$ws.Range("K94").Value = 17100
$ws.Range("L94").Value = 26000.001
$ws.Range("M94").Value = -16424
$ws.Range("N94").Value = -27352.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 24000
$ws.Range("J58").Value = 24000
$ws.Range("L58").Value = 24000
$ws.Range("N58").Value = -24554
$ws.Range("H97").Value = 1039
$ws.Range("I97").Value = 888.75
$ws.Range("K97").Value = 888.75
$ws.Range("M97").Value = -392.75
$ws.Range("H132").Value = 3007
$ws.Range("J132").Value = 2996.25
$ws.Range("L132").Value = 8988.75
$ws.Range("N132").Value = -14048.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2756.303
$ws.Range("J22").Value = 3482.5334
$ws.Range("L22").Value = 3482.5334
$ws.Range("N22").Value = -4072.5334
$ws.Range("H27").Value = 2756.303
$ws.Range("J27").Value = 3482.5334
$ws.Range("L27").Value = 3482.5334
$ws.Range("N27").Value = -3696.5334
$ws.Range("H46").Value = 3089.5715
$ws.Range("I46").Value = 2237.1428
$ws.Range("J46").Value = 3942
$ws.Range("K46").Value = 2237.1428
$ws.Range("L46").Value = 3942
$ws.Range("M46").Value = -2049.1428
$ws.Range("N46").Value = -4318
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 5000
$ws.Range("K56").Value = 5000
$ws.Range("M56").Value = -4309
$ws.Range("H132").Value = 2479.5938
$ws.Range("I132").Value = 2230.8572
$ws.Range("J132").Value = 4220.75
$ws.Range("K132").Value = 6692.571599999999
$ws.Range("L132").Value = 12662.25
$ws.Range("M132").Value = -4162.571599999999
$ws.Range("N132").Value = -17722.25
$ws.Range("H136").Value = 6650
$ws.Range("I136").Value = 6400
$ws.Range("K136").Value = 19200
$ws.Range("M136").Value = -16650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H132").Value = 1425.88
$ws.Range("I132").Value = 1308.2084
$ws.Range("K132").Value = 3924.6252
$ws.Range("M132").Value = -1394.6252
